# Apply the numeric cell updates from the scheduled-runner diff.
# (H:N "profit" columns across the 8 item-type sheets.)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1632.6111
$ws.Range("I19").Value = 1250.125
$ws.Range("J19").Value = 1938.6
$ws.Range("K19").Value = 1250.125
$ws.Range("L19").Value = 1938.6
$ws.Range("M19").Value = -1075.125
$ws.Range("N19").Value = -2288.6
$ws.Range("H33").Value = 859.23334
$ws.Range("I33").Value = 234.54167
$ws.Range("J33").Value = 3358
$ws.Range("K33").Value = 234.54167
$ws.Range("L33").Value = 3358
$ws.Range("M33").Value = -5.541670000000011
$ws.Range("N33").Value = -3816
$ws.Range("H98").Value = 7666.1113
$ws.Range("I98").Value = 7572.143
$ws.Range("J98").Value = 7995
$ws.Range("K98").Value = 7572.143
$ws.Range("L98").Value = 7995
$ws.Range("M98").Value = -6074.143
$ws.Range("N98").Value = -10991
$ws.Range("H112").Value = 1062.421
$ws.Range("I112").Value = 759.2857
$ws.Range("J112").Value = 1239.25
$ws.Range("K112").Value = 2277.8571
$ws.Range("L112").Value = 3717.75
$ws.Range("M112").Value = -1169.8571
$ws.Range("N112").Value = -5933.75
$ws.Range("H113").Value = 101814.4
$ws.Range("I113").Value = 501102.5
$ws.Range("J113").Value = 1992.375
$ws.Range("K113").Value = 501102.5
$ws.Range("L113").Value = 1992.375
$ws.Range("M113").Value = -497848.5
$ws.Range("N113").Value = -8500.375
$ws.Range("H122").Value = 7666.1113
$ws.Range("I122").Value = 7572.143
$ws.Range("J122").Value = 7995
$ws.Range("K122").Value = 22716.429
$ws.Range("L122").Value = 23985
$ws.Range("M122").Value = -20266.429
$ws.Range("N122").Value = -28885
$ws.Range("H127").Value = 26317478
$ws.Range("I127").Value = 447.1
$ws.Range("J127").Value = 35716416
$ws.Range("K127").Value = 1341.3
$ws.Range("L127").Value = 107149248
$ws.Range("M127").Value = 3618.7
$ws.Range("N127").Value = -107159168
$ws.Range("H129").Value = 2272.2134
$ws.Range("I129").Value = 6050.222
$ws.Range("J129").Value = 1079.1578
$ws.Range("K129").Value = 18150.666
$ws.Range("L129").Value = 3237.4734
$ws.Range("M129").Value = -13150.666
$ws.Range("N129").Value = -13237.4734
$ws.Range("H132").Value = 5324012
$ws.Range("I132").Value = 5957746.5
$ws.Range("J132").Value = 642.2
$ws.Range("K132").Value = 17873239.5
$ws.Range("L132").Value = 1926.6
$ws.Range("M132").Value = -17870709.5
$ws.Range("N132").Value = -6986.6
$ws.Range("H135").Value = 1001.04254
$ws.Range("I135").Value = 643.1905
$ws.Range("J135").Value = 4007
$ws.Range("K135").Value = 5788.7145
$ws.Range("L135").Value = 36063
$ws.Range("M135").Value = -3253.7145
$ws.Range("N135").Value = -41133
$ws.Range("H137").Value = 1011.51666
$ws.Range("I137").Value = 1024.5106
$ws.Range("J137").Value = 964.53845
$ws.Range("K137").Value = 3073.5318
$ws.Range("L137").Value = 2893.61535
$ws.Range("M137").Value = -523.5318000000002
$ws.Range("N137").Value = -7993.61535
$ws.Range("H138").Value = 1442.7037
$ws.Range("I138").Value = 1001.2093
$ws.Range("J138").Value = 3168.5454
$ws.Range("K138").Value = 3003.6279
$ws.Range("L138").Value = 9505.636200000001
$ws.Range("M138").Value = 2136.3721
$ws.Range("N138").Value = -19785.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 599.34
$ws.Range("I32").Value = 575.53845
$ws.Range("J32").Value = 840
$ws.Range("K32").Value = 575.53845
$ws.Range("L32").Value = 840
$ws.Range("M32").Value = -288.53845
$ws.Range("N32").Value = -1414
$ws.Range("H39").Value = 13204
$ws.Range("I39").Value = 5008
$ws.Range("J39").Value = 21400
$ws.Range("K39").Value = 5008
$ws.Range("L39").Value = 21400
$ws.Range("M39").Value = -4488
$ws.Range("N39").Value = -22440
$ws.Range("H42").Value = 14765.6
$ws.Range("I42").Value = 11028
$ws.Range("J42").Value = 15700
$ws.Range("K42").Value = 11028
$ws.Range("L42").Value = 15700
$ws.Range("M42").Value = -10542
$ws.Range("N42").Value = -16672
$ws.Range("H61").Value = 1158.2433
$ws.Range("I61").Value = 957.5294
$ws.Range("J61").Value = 3433
$ws.Range("K61").Value = 957.5294
$ws.Range("L61").Value = 3433
$ws.Range("M61").Value = -745.5294
$ws.Range("N61").Value = -3857
$ws.Range("H74").Value = 604.53845
$ws.Range("I74").Value = 458.1875
$ws.Range("J74").Value = 838.7
$ws.Range("K74").Value = 458.1875
$ws.Range("L74").Value = 838.7
$ws.Range("M74").Value = 415.8125
$ws.Range("N74").Value = -2586.7
$ws.Range("H77").Value = 604.53845
$ws.Range("I77").Value = 458.1875
$ws.Range("J77").Value = 838.7
$ws.Range("K77").Value = 2290.9375
$ws.Range("L77").Value = 4193.5
$ws.Range("M77").Value = 2077.0625
$ws.Range("N77").Value = -12929.5
$ws.Range("H122").Value = 1625.625
$ws.Range("I122").Value = 1665.1666
$ws.Range("J122").Value = 1507
$ws.Range("K122").Value = 4995.4998
$ws.Range("L122").Value = 4521
$ws.Range("M122").Value = -2545.4998
$ws.Range("N122").Value = -9421
$ws.Range("H132").Value = 2244.5952
$ws.Range("I132").Value = 2282.3896
$ws.Range("J132").Value = 1828.8572
$ws.Range("K132").Value = 6847.168799999999
$ws.Range("L132").Value = 5486.571599999999
$ws.Range("M132").Value = -4317.168799999999
$ws.Range("N132").Value = -10546.5716
$ws.Range("H136").Value = 1158.2433
$ws.Range("I136").Value = 957.5294
$ws.Range("J136").Value = 3433
$ws.Range("K136").Value = 2872.5882
$ws.Range("L136").Value = 10299
$ws.Range("M136").Value = -322.5882000000001
$ws.Range("N136").Value = -15399

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 39800
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 39800
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 39800
$ws.Range("N9").Value = -40136
$ws.Range("H134").Value = 2699.1052
$ws.Range("I134").Value = 2387.6572
$ws.Range("J134").Value = 6332.6665
$ws.Range("K134").Value = 7162.971600000001
$ws.Range("L134").Value = 18997.9995
$ws.Range("M134").Value = -4627.971600000001
$ws.Range("N134").Value = -24067.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18787.723
$ws.Range("I31").Value = 1118.0741
$ws.Range("J31").Value = 51689.83
$ws.Range("K31").Value = 1118.0741
$ws.Range("L31").Value = 51689.83
$ws.Range("M31").Value = -823.0741
$ws.Range("N31").Value = -52279.83
$ws.Range("H34").Value = 18787.723
$ws.Range("I34").Value = 1118.0741
$ws.Range("J34").Value = 51689.83
$ws.Range("K34").Value = 1118.0741
$ws.Range("L34").Value = 51689.83
$ws.Range("M34").Value = -916.0741
$ws.Range("N34").Value = -52093.83
$ws.Range("H35").Value = 5875
$ws.Range("I35").Value = 2812.5
$ws.Range("J35").Value = 12000
$ws.Range("K35").Value = 2812.5
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = -2518.5
$ws.Range("N35").Value = -12588
$ws.Range("H132").Value = 2575.0544
$ws.Range("I132").Value = 2281.8333
$ws.Range("J132").Value = 4585.7144
$ws.Range("K132").Value = 6845.499899999999
$ws.Range("L132").Value = 13757.1432
$ws.Range("M132").Value = -4315.499899999999
$ws.Range("N132").Value = -18817.1432
$ws.Range("H134").Value = 907.62
$ws.Range("I134").Value = 887.6047
$ws.Range("J134").Value = 1030.5714
$ws.Range("K134").Value = 2662.8141
$ws.Range("L134").Value = 3091.7142
$ws.Range("M134").Value = -127.8141000000001
$ws.Range("N134").Value = -8161.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 736469.7
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 736469.7
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 2209409.1
$ws.Range("N37").Value = -2209633.1
$ws.Range("H55").Value = 13182.611
$ws.Range("I55").Value = 25455
$ws.Range("J55").Value = 9676.214
$ws.Range("K55").Value = 76365
$ws.Range("L55").Value = 29028.642
$ws.Range("M55").Value = -76188
$ws.Range("N55").Value = -29382.642
$ws.Range("H131").Value = 1352.8721
$ws.Range("I131").Value = 826.125
$ws.Range("J131").Value = 1406.8975
$ws.Range("K131").Value = 2478.375
$ws.Range("L131").Value = 4220.6925
$ws.Range("M131").Value = 2561.625
$ws.Range("N131").Value = -14300.6925

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 10888.444
$ws.Range("I46").Value = 6000
$ws.Range("J46").Value = 11499.5
$ws.Range("K46").Value = 6000
$ws.Range("L46").Value = 11499.5
$ws.Range("M46").Value = -5844
$ws.Range("N46").Value = -11811.5
$ws.Range("H57").Value = 13266.25
$ws.Range("I57").Value = 12721.667
$ws.Range("J57").Value = 14900
$ws.Range("K57").Value = 12721.667
$ws.Range("L57").Value = 14900
$ws.Range("M57").Value = -11901.667
$ws.Range("N57").Value = -16540
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1704.2609
$ws.Range("I132").Value = 1645.3636
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4936.0908
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2406.0908
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2633.6667
$ws.Range("I122").Value = 2560.4
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 7681.200000000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -5231.200000000001
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 2243.7083
$ws.Range("I132").Value = 2248.8445
$ws.Range("J132").Value = 2166.6667
$ws.Range("K132").Value = 6746.5335
$ws.Range("L132").Value = 6500.000100000001
$ws.Range("M132").Value = -4216.5335
$ws.Range("N132").Value = -11560.0001
$ws.Range("H136").Value = 1749.909
$ws.Range("I136").Value = 1644.3334
$ws.Range("J136").Value = 2225
$ws.Range("K136").Value = 4933.0002
$ws.Range("L136").Value = 6675
$ws.Range("M136").Value = -2383.0002
$ws.Range("N136").Value = -11775

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 142866850
$ws.Range("I39").Value = 250004690
$ws.Range("J39").Value = 16400
$ws.Range("K39").Value = 250004690
$ws.Range("L39").Value = 16400
$ws.Range("M39").Value = -250004277
$ws.Range("N39").Value = -17226
$ws.Range("H132").Value = 2041.5283
$ws.Range("I132").Value = 2062.2654
$ws.Range("J132").Value = 1787.5
$ws.Range("K132").Value = 6186.796200000001
$ws.Range("L132").Value = 5362.5
$ws.Range("M132").Value = -3656.796200000001
$ws.Range("N132").Value = -10422.5
$ws.Range("H136").Value = 475.5849
$ws.Range("I136").Value = 314.22
$ws.Range("J136").Value = 3165
$ws.Range("K136").Value = 942.6600000000001
$ws.Range("L136").Value = 9495
$ws.Range("M136").Value = 1607.34
$ws.Range("N136").Value = -14595
